$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-6), column A is trial number
$data = @(
    @(1, 2, 9, 7, 8, 5, -1, 12, 5),
    @(2, 0, 9, 4, 7, 4, -2, 23, 5),
    @(3, 4, 5, 5, 0, 1, -5, 56, 5),
    @(4, 3, 8, 6, 5, 3, -3, 34, 5),
    @(5, 1, 6, 3, 2, 2, -4, 45, 5)
)

$row = 2
foreach ($rowValues in $data) {
    $col = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
    $row++
}
